# Weekly fruit/vegetable price update: a new daily record is inserted at the
# top of the dated price table (row 373), pushing all existing records for
# rows 373-502 down by one row (to 374-503). The sheet's used range grows
# from A1:R502 to A1:R503.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the first data row of the block, shifting
# the existing rows (and all their formatting) down by one.
$ws.Rows("373:373").Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(373, 1).Value = 10
$ws.Cells.Item(373, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(373, 3).Value = "La Araucanía"
$ws.Cells.Item(373, 4).Value = 45120
$ws.Cells.Item(373, 5).Value = 9
$ws.Cells.Item(373, 6).Value = 100112044
$ws.Cells.Item(373, 7).Value = "Perejil"
$ws.Cells.Item(373, 8).Value = "Sin especificar"
$ws.Cells.Item(373, 9).Value = "Primera"
$ws.Cells.Item(373, 10).Value = 100
$ws.Cells.Item(373, 11).Value = 4000
$ws.Cells.Item(373, 12).Value = 4000
$ws.Cells.Item(373, 13).Value = 4000
$ws.Cells.Item(373, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(373, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(373, 16).Value = 1333
$ws.Cells.Item(373, 17).Value = 3
$ws.Cells.Item(373, 18).Value = "Hortaliza"
